$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($addr, $val)
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '62.556.68'
Set-TextValue 'E2' '  +2.03%  '
Set-TextValue 'D3' '3.442.21'
Set-TextValue 'E3' '  +3.05%  '
Set-TextValue 'E4' '  -0.04%  '
Set-TextValue 'D5' '406.96'
Set-TextValue 'E5' '  +2.24%  '
Set-TextValue 'D6' '131.35'
Set-TextValue 'E6' '  +5.07%  '
Set-TextValue 'D7' '0.602'
Set-TextValue 'E7' '  +2.85%  '
Set-TextValue 'D8' '0.999'
Set-TextValue 'E8' '  -0.06%  '
Set-TextValue 'D9' '0.700'
Set-TextValue 'E9' '  +7.01%  '
Set-TextValue 'D10' '0.141'
Set-TextValue 'E10' '  +19.63%  '
Set-TextValue 'D11' '42.50'
Set-TextValue 'E11' '  +5.03%  '
Set-TextValue 'D13' '8.64'
Set-TextValue 'E13' '  +5.23%  '
Set-TextValue 'D14' '20.00'
Set-TextValue 'E14' '  +4.39%  '
Set-TextValue 'D15' '3.431.50'
Set-TextValue 'E15' '  +5.01%  '
Set-TextValue 'D16' '62.583.63'
Set-TextValue 'E16' '  +2.18%  '
Set-TextValue 'D17' '11.61'
Set-TextValue 'E17' '  +4.65%  '
Set-TextValue 'E18' '  +2.55%  '
Set-TextValue 'D19' '0.0000164'
Set-TextValue 'E19' '  +31.29%  '
Set-TextValue 'E20' '  +0.67%  '
Set-TextValue 'D21' '84.77'
Set-TextValue 'E21' '  +6.68%  '
Set-TextValue 'D22' '316.72'
Set-TextValue 'E22' '  +6.63%  '
Set-TextValue 'D23' '12.93'
Set-TextValue 'E23' '  +2.28%  '
Set-TextValue 'D24' '3.19'
Set-TextValue 'E24' '  +3.93%  '
Set-TextValue 'D25' '4.76'
Set-TextValue 'E25' '  +0.14%  '
Set-TextValue 'D26' '30.11'
Set-TextValue 'E26' '  +4.62%  '
Set-TextValue 'D27' '8.20'
Set-TextValue 'E27' '  +1.36%  '
Set-TextValue 'D28' '7.84'
Set-TextValue 'E28' '  +6.12%  '
Set-TextValue 'D29' '2.78'
Set-TextValue 'E29' '  +11.39%  '
Set-TextValue 'B30' 'InjectiveProtocol'
Set-TextValue 'C30' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D30' '44.75'
Set-TextValue 'E30' '  +9.65%  '
Set-TextValue 'B31' 'Kaspa'
Set-TextValue 'C31' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D31' '0.174'
Set-TextValue 'E31' '  +1.84%  '
Set-TextValue 'D32' '0.116'
Set-TextValue 'E32' '  +3.09%  '
Set-TextValue 'D33' '11.52'
Set-TextValue 'E33' '  +2.31%  '
Set-TextValue 'D34' '1.00'
Set-TextValue 'D35' '0.0489'
Set-TextValue 'E35' '  +3.18%  '
Set-TextValue 'D36' '51.44'
Set-TextValue 'E36' '  -0.72%  '
Set-TextValue 'E37' '  +0.10%  '
Set-TextValue 'D38' '2.98'
Set-TextValue 'E38' '  +3.41%  '
Set-TextValue 'D39' '3.35'
Set-TextValue 'E39' '  +0.02%  '
Set-TextValue 'E40' '  +15.29%  '
Set-TextValue 'D41' '143.89'
Set-TextValue 'E41' '  +5.13%  '
Set-TextValue 'D42' '0.127'
Set-TextValue 'E42' '  +4.70%  '
Set-TextValue 'E43' '  +2.59%  '
Set-TextValue 'D44' '17.02'
Set-TextValue 'E44' '  +3.44%  '
Set-TextValue 'D45' '3.95'
Set-TextValue 'E45' '  +2.89%  '
Set-TextValue 'E46' '  +0.32%  '
Set-TextValue 'D47' '21.55'
Set-TextValue 'E47' '  +3.26%  '
Set-TextValue 'D48' '2.110.48'
Set-TextValue 'E48' '  +1.15%  '
Set-TextValue 'E49' '  +11.10%  '
Set-TextValue 'D50' '2.33'
Set-TextValue 'E50' '  +2.01%  '
Set-TextValue 'D51' '1.08'
Set-TextValue 'E51' '  +31.20%  '
